$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 20: equipment necrot3 body
$ws.Range("B20").Value = "itd_body_necrot3"
$ws.Range("A20").Value = "it_eq_body_necrot3"
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = "100, 200"
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 30
$ws.Range("K20").Value = 10
$ws.Range("L20").Value = 10
$ws.Range("M20").Value = 20
$ws.Range("O20").Value = 30
$ws.Range("P20").Value = 5
$ws.Range("Q20").Value = 5
$ws.Range("R20").Value = 5
$ws.Range("S20").Value = 5
$ws.Range("Y20").Value = "res/assets/equipment/body/spritesheet_body_necrot3.png"

# Move/refresh the active selection, matching the post-edit view state
[void]$ws.Range("AB21").Select()
